$wb = $excel.ActiveWorkbook

$oldTimestamp = "February 03 2026 17.29.55 EST"
$newTimestamp = "February 03 2026 18.05.36 EST"

foreach ($ws in $wb.Worksheets) {
    [void]$ws.Cells.Replace($oldTimestamp, $newTimestamp)
}
